$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'youth knee sleeves with pads'
    'youth gloves and knee pads'
    'mens tights 3 4'
    'black leggings capri'
    'knee pad tights'
    'knee leggings'
    'mens compression pants green'
    'men tights leggings'
    'knee pad workout'
    'hex knee pads basketball'
    'white knee pads basketball'
    'knee pad basketball leggings'
    'capri compression tights'
    'compression workout capri'
    'knee pad leggings youth'
    'all white knee pads for basketball'
    'basketball compression pants with knee pads youth'
    'triple 8 knee pads youth'
    'mens leggings compression under armour'
    'dirt bike knee pads youth'
    'nike leggings men basketball'
    'gym shark mens leggings'
    'multicam pants with knee pads'
    'tactical pants with knee pads for men'
    'excersize gear for men'
    'men capri leggings'
    'mens leggings black'
    'under armour tights'
    'guys tights'
    'leggings for basketball'
    'dri fit compression pants men'
    'pink compression pants men'
    'men workout compression pants'
    'mens football gear'
    'compression pants with pads basketball'
    'mens white leggings compression'
    'cold gear mens'
    'basketball leggings with kneepads'
    'mens athletic leggins'
    'mens compression tights basketball'
    'under armour cold gear compression pants men'
    'green mens compression leggings'
    'mens compression tights leggings'
    'compression knee pads for basketball'
    'legging basketball men'
    'basketball compression pants with padded knees'
    'knee pad pants'
    'kneepad honeycomb'
    'black basketball knee pads'
    'cheap knee pads for basketball'
    'compression pants'
    'knee pads pants'
    'basketball youth compression pants'
    'knee pad hex'
    'men tights and leggings'
    'knee pad lacrosse'
    'mens compression tights'
    'patella guard'
    'squat pad knee'
    'best knee pads basketball'
    'kneepad basketball'
    'compression mens running pants'
    'knee protection soccer'
    'sliding leg guard'
    'athletic leggings youth'
    'basketball tights boys'
    'knee pads xxl'
    'wrestling clothes for boys'
    'boys xxl baseball pants'
    'knee pads girls volleyball'
    'youth boys leggings sports'
    'mens compression running pants'
    'black knee pads wrestling'
    'black knee pads youth'
    'hex pad'
    'men spandex pants'
    'knee protectors for men'
    'volleyball gear'
    'capris for men'
    'youth black baseball pants'
    'knee pads professional'
    'mens mesh pants'
    'knee pad work pants'
    'knee pads small'
    'youth leggings sports'
    'fitness knee pads'
    'knee pad sports'
    'basketball clothes'
    'girls workout leggings'
    'baseball pants mens'
    'kids basketball leggings with knee pads'
    'youth asics wrestling knee pads'
    'tesla wintergear for men'
    'mcdavid basketball knee pads white'
    'thermo ball mens'
    'mens nike thermal training pants'
    'man winter leggings'
    'underware pants men'
    'under armor youth basketball compression pants'
    'black capris'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

